$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert two new columns before D (one for Q4'18, one for Q3'18) ---
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number/style formats from column F (which used to be D, holding the
# "date"/"value" styles) into the two freshly inserted D:E columns, only for
# the row ranges that actually carry data (skip the bare section-title rows).
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)

$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)

$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Write the updated figures (new quarters + restated prior figures) ---
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 1749000
$ws.Cells.Item(8, 5).Value = 1747700
$ws.Cells.Item(8, 6).Value = 3583400
$ws.Cells.Item(9, 4).Value = 1494900
$ws.Cells.Item(9, 5).Value = 1512800
$ws.Cells.Item(9, 6).Value = 3078700
$ws.Cells.Item(10, 4).Value = 254100
$ws.Cells.Item(10, 5).Value = 234900
$ws.Cells.Item(10, 6).Value = 504700
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 700
$ws.Cells.Item(14, 5).Value = 900
$ws.Cells.Item(14, 6).Value = 2100
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(17, 4).Value = 1690400
$ws.Cells.Item(17, 5).Value = 1697800
$ws.Cells.Item(17, 6).Value = 3458500
$ws.Cells.Item(18, 4).Value = 58600
$ws.Cells.Item(18, 5).Value = 49900
$ws.Cells.Item(18, 6).Value = 124900
$ws.Cells.Item(20, 4).Value = 1700
$ws.Cells.Item(20, 5).Value = -600
$ws.Cells.Item(20, 6).Value = 200
$ws.Cells.Item(21, 4).Value = 70300
$ws.Cells.Item(21, 5).Value = 58800
$ws.Cells.Item(21, 6).Value = 143000
$ws.Cells.Item(22, 4).Value = 5600
$ws.Cells.Item(22, 5).Value = 6100
$ws.Cells.Item(22, 6).Value = 11100
$ws.Cells.Item(23, 4).Value = 54700
$ws.Cells.Item(23, 5).Value = 43200
$ws.Cells.Item(23, 6).Value = 114000
$ws.Cells.Item(24, 4).Value = 13300
$ws.Cells.Item(24, 5).Value = 11100
$ws.Cells.Item(24, 6).Value = 29500
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 41400
$ws.Cells.Item(26, 5).Value = 32200
$ws.Cells.Item(26, 6).Value = 84500
$ws.Cells.Item(27, 4).Value = 41400
$ws.Cells.Item(27, 5).Value = 32200
$ws.Cells.Item(27, 6).Value = 84500
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = 5600
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = -1700
$ws.Cells.Item(32, 5).Value = 600
$ws.Cells.Item(32, 6).Value = -200
$ws.Cells.Item(33, 4).Value = 47000
$ws.Cells.Item(33, 5).Value = 32200
$ws.Cells.Item(33, 6).Value = 84500
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 47000
$ws.Cells.Item(35, 5).Value = 32200
$ws.Cells.Item(35, 6).Value = 84500
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 142700
$ws.Cells.Item(41, 5).Value = 111100
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(43, 4).Value = 1931700
$ws.Cells.Item(43, 5).Value = 1682000
$ws.Cells.Item(44, 4).Value = 148500
$ws.Cells.Item(44, 5).Value = 171800
$ws.Cells.Item(45, 4).Value = 115700
$ws.Cells.Item(45, 5).Value = 103800
$ws.Cells.Item(46, 4).Value = 2338600
$ws.Cells.Item(46, 5).Value = 2068700
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(48, 4).Value = 73000
$ws.Cells.Item(48, 5).Value = 74100
$ws.Cells.Item(49, 4).Value = 279000
$ws.Cells.Item(49, 5).Value = 283700
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 85400
$ws.Cells.Item(52, 5).Value = 84100
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 2775900
$ws.Cells.Item(54, 5).Value = 2510500
$ws.Cells.Item(57, 4).Value = 1282200
$ws.Cells.Item(57, 5).Value = 995600
$ws.Cells.Item(58, 4).Value = 1400
$ws.Cells.Item(58, 5).Value = 17400
$ws.Cells.Item(59, 4).Value = 253000
$ws.Cells.Item(59, 5).Value = 243800
$ws.Cells.Item(60, 4).Value = 1536700
$ws.Cells.Item(60, 5).Value = 1256700
$ws.Cells.Item(61, 4).Value = 195500
$ws.Cells.Item(61, 5).Value = 251300
$ws.Cells.Item(62, 4).Value = 56800
$ws.Cells.Item(62, 5).Value = 59400
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 1789000
$ws.Cells.Item(66, 5).Value = 1567500
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 704700
$ws.Cells.Item(72, 5).Value = 657600
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 987000
$ws.Cells.Item(76, 5).Value = 943000
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 47000
$ws.Cells.Item(81, 5).Value = 32200
$ws.Cells.Item(81, 6).Value = 84500
$ws.Cells.Item(83, 4).Value = 10000
$ws.Cells.Item(83, 5).Value = 9500
$ws.Cells.Item(83, 6).Value = 17900
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 45400
$ws.Cells.Item(89, 5).Value = -103700
$ws.Cells.Item(89, 6).Value = 351000
$ws.Cells.Item(91, 4).Value = -4200
$ws.Cells.Item(91, 5).Value = -2400
$ws.Cells.Item(91, 6).Value = -10600
$ws.Cells.Item(91, 9).Value = -5600
$ws.Cells.Item(91, 10).Value = -200
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -4200
$ws.Cells.Item(94, 5).Value = -76900
$ws.Cells.Item(94, 6).Value = -10600
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(96, 5).Value = 0
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = -6900
$ws.Cells.Item(100, 5).Value = 40400
$ws.Cells.Item(100, 6).Value = -192500
$ws.Cells.Item(101, 4).Value = -2600
$ws.Cells.Item(101, 5).Value = 3100
$ws.Cells.Item(101, 6).Value = -5500
$ws.Cells.Item(101, 8).Value = -3400
$ws.Cells.Item(102, 4).Value = 31600
$ws.Cells.Item(102, 5).Value = -137100
$ws.Cells.Item(102, 6).Value = 142300

# Re-fit the now-wider used range (matches the author's "bestFit" columns)
$ws.Columns.Item("A:M").AutoFit()

Write-Output "done"
